$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.690834164619446
$ws.Range("B1").Value = 1.922191023826599
$ws.Range("C1").Value = 5.188684463500977
$ws.Range("D1").Value = 1.332329511642456
$ws.Range("E1").Value = 0.7423061728477478
